$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update company names (rows 3 and 4 swapped identities)
$ws.Range('B3').Value = 'Kuznetsky Bank Public Joint Stock Company (MISX:KUZB)'
$ws.Range('B4').Value = 'The Russian Public Joint-Stock Commercial Roads Bank (Public joint-stock company) (MISX:RDRB)'

# Row 2 numeric updates
$ws.Range('D2').Value = 0.06632
$ws.Range('E2').Value = 0.7044999999999999
$ws.Range('K2').Value = 2.146
$ws.Range('L2').Value = 0.09495575221238939
$ws.Range('M2').Value = 0.193
$ws.Range('N2').Value = 0.003893483962073835
$ws.Range('O2').Value = 0.08993476234855546
$ws.Range('P2').Value = 0.193
$ws.Range('Q2').Value = 0.003893483962073835
$ws.Range('R2').Value = 0.08993476234855546
$ws.Range('U2').Value = 51.26000000000001
$ws.Range('V2').Value = 1.034093201533185
$ws.Range('W2').Value = 0.04650849944967592
$ws.Range('X2').Value = 0.05366381183770717
$ws.Range('Y2').Value = -0.00715531238803125
$ws.Range('Z2').Value = 0.6764036872979767
$ws.Range('AB2').Value = 0.0505699239683618
$ws.Range('AC2').Value = -0.0505699239683618
$ws.Range('AD2').Value = 12.369
$ws.Range('AF2').Value = 12.369
$ws.Range('AG2').Value = -38.89100000000001
$ws.Range('AH2').Value = 0.1996964755646685
$ws.Range('AI2').Value = 0.2038431747392014
$ws.Range('AJ2').Value = -3.641820395168089
$ws.Range('AK2').Value = -4.128994585412466

# Row 3 numeric updates
$ws.Range('D3').Value = 0.00664
$ws.Range('E3').Value = -0.366
$ws.Range('K3').Value = 0.046
$ws.Range('L3').Value = 0.008846153846153846
$ws.Range('M3').Value = 0.193
$ws.Range('N3').Value = 0.04742014742014742
$ws.Range('O3').Value = 4.195652173913044
$ws.Range('P3').Value = 0.193
$ws.Range('Q3').Value = 0.04742014742014742
$ws.Range('R3').Value = 4.195652173913044
$ws.Range('U3').Value = 6.06
$ws.Range('V3').Value = 1.488943488943489
$ws.Range('W3').Value = 0.004781704781704782
$ws.Range('X3').Value = 0.0534787655650044
$ws.Range('Y3').Value = -0.04869706078329962
$ws.Range('Z3').Value = 0.5990783410138251
$ws.Range('AB3').Value = 0.05015651416672787
$ws.Range('AC3').Value = -0.05015651416672787
$ws.Range('AD3').Value = 0.969
$ws.Range('AF3').Value = 0.969
$ws.Range('AG3').Value = -5.090999999999999
$ws.Range('AH3').Value = 0.1923000595356221
$ws.Range('AI3').Value = 0.1103770361088962
$ws.Range('AJ3').Value = 4.986287952987271
$ws.Range('AK3').Value = -1.872379551305626

# Row 4 numeric updates
$ws.Range('D4').Value = 0.126
$ws.Range('E4').Value = 1.775
$ws.Range('K4').Value = 2.1
$ws.Range('L4').Value = 0.1206896551724138
$ws.Range('M4').Value = -0
$ws.Range('N4').Value = -0
$ws.Range('O4').Value = -0
$ws.Range('P4').Value = -0
$ws.Range('Q4').Value = -0
$ws.Range('R4').Value = -0
$ws.Range('U4').Value = 45.2
$ws.Range('V4').Value = 0.9934065934065934
$ws.Range('W4').Value = 0.08823529411764706
$ws.Range('X4').Value = 0.05384885811040994
$ws.Range('Y4').Value = 0.03438643600723712
$ws.Range('Z4').Value = 0.7035419699175157
$ws.Range('AB4').Value = 0.05098333376999572
$ws.Range('AC4').Value = -0.05098333376999572
$ws.Range('AD4').Value = 11.4
$ws.Range('AF4').Value = 11.4
$ws.Range('AG4').Value = -33.8
$ws.Range('AH4').Value = 0.2003514938488576
$ws.Range('AI4').Value = 0.2196531791907514
$ws.Range('AJ4').Value = -2.88888888888889
$ws.Range('AK4').Value = -5.044776119402989

# T4 cell removed entirely in the updated data
$ws.Range('T4').ClearContents()
